# ---------------------------------------------------------------------------
# EPushButton (弹出按钮) event-function properties
#   - QStudioSCADA / QSCADARunTime: add the "事件功能" (Event Function) sheet
#     listing every new event action, and refresh the dev-progress status of
#     the two screen-editor rows that already shipped this (弹出按钮, 输入编辑框).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# Colors (VBA/COM "Color" is BGR-encoded, not RGB)
$colGreen  = 5287936   # RGB 00B050 -> 已完成 (Completed)
$colYellow = 65535     # RGB FFFF00 -> 进行中 (In progress)
$colRed    = 255       # RGB FF0000 -> 未开始 (Not started)
$colPurple = 10642560  # RGB 8064A2 (theme accent4) -> table grid-line color

# ---------------------------------------------------------------------------
# 1. Sheet1 "画面编辑器" - bump the status of the two controls whose event
#    properties just landed.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# 弹出按钮 (popup button): 进行中 -> 已完成
$ws1.Range("B10").Value = "已完成"
$ws1.Range("B10").Interior.Color = $colGreen

# 输入编辑框 (input edit box): 未开始 -> 进行中
$ws1.Range("B11").Value = "进行中"
$ws1.Range("B11").Interior.Color = $colYellow

# ---------------------------------------------------------------------------
# 2. Sheet2 gets renamed to "事件功能" (Event Function) and becomes the sheet
#    tracking every new EPushButton event action - all starting "未开始".
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "事件功能"

$eventFunctions = @(
    "切换画面", "返回画面", "设置时间", "执行脚本", "设置变量值",
    "获取变量的值", "切换变量状态", "设置系统变量值", "获取系统变量值",
    "变量值自增", "变量值自减", "隐藏控件", "显示控件", "失效控件",
    "生效控件", "偏移控件", "闪烁控件", "停止闪烁"
)

$lastRow = 43
$lastDataRow = $eventFunctions.Count   # 18

# Column headers / widths matching the 功能模块-状态 layout from sheet1.
$ws2.Columns.Item(1).ColumnWidth = 27.5
$ws2.Columns.Item(2).ColumnWidth = 25.36

for ($i = 0; $i -lt $eventFunctions.Count; $i++) {
    $r = $i + 1
    $ws2.Range("A$r").Value = $eventFunctions[$i]
    $ws2.Range("A$r").HorizontalAlignment = -4108   # xlCenter
    $ws2.Range("B$r").Value = "未开始"
    $ws2.Range("B$r").Interior.Color = $colRed
}

# Thin purple-ish grid around the A:B table, row by row (top edge on every
# row + left edge on every A cell + right edge on every B cell), with the
# bottom edge closing the table only on the very last row.
for ($r = 1; $r -le $lastRow; $r++) {
    $aCell = $ws2.Range("A$r")
    $bCell = $ws2.Range("B$r")

    $aCell.Borders.Item(8).LineStyle = 1      # xlEdgeTop
    $aCell.Borders.Item(8).Color = $colPurple
    $aCell.Borders.Item(7).LineStyle = 1      # xlEdgeLeft
    $aCell.Borders.Item(7).Color = $colPurple

    $bCell.Borders.Item(8).LineStyle = 1      # xlEdgeTop
    $bCell.Borders.Item(8).Color = $colPurple
    $bCell.Borders.Item(10).LineStyle = 1     # xlEdgeRight
    $bCell.Borders.Item(10).Color = $colPurple
}
$ws2.Range("A$lastRow").Borders.Item(9).LineStyle = 1   # xlEdgeBottom
$ws2.Range("A$lastRow").Borders.Item(9).Color = $colPurple
$ws2.Range("B$lastRow").Borders.Item(9).LineStyle = 1
$ws2.Range("B$lastRow").Borders.Item(9).Color = $colPurple

# Status dropdown across the whole table body, same list as sheet1's table.
$null = $ws2.Range("B1:B$lastRow").Validation.Add(3, 1, 1, '"未开始,已完成,进行中"')

# ---------------------------------------------------------------------------
# 3. Selection / active-tab bookkeeping - the edit session ends with the new
#    "事件功能" sheet active (and scrolled to C20), while sheet1 is left
#    selected over its (now taller) data range.
# ---------------------------------------------------------------------------
$ws1.Activate()
$null = $ws1.Range("A2:B44").Select()

$ws2.Activate()
$null = $ws2.Range("C20").Select()
